$d = $word.ActiveDocument

$startPos = $d.Paragraphs(1).Range.Start
$endPos = $d.Paragraphs(6).Range.Start
$r = $d.Range($startPos, $endPos)
$r.Delete()

Write-Output ("ParaCount=" + $d.Paragraphs.Count)
Write-Output ("P1=" + $d.Paragraphs(1).Range.Text)
Write-Output ("P2=" + $d.Paragraphs(2).Range.Text)
